$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Productivity %"
$ws.Range("E1").Value = "Quality %"
$ws.Range("F1").Value = "Present %"
$ws.Range("G1").Value = "Final %"
$ws.Range("H1").Value = "Grade"
$ws.Range("I1").Value = "Absent Days"
$ws.Range("J1").Value = "HR Comments"
$ws.Range("K1").ClearContents()

$ws.Range("D1:J1").Select()
